# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted at the top of the data block
# (row 152), pushing all the existing records down by one row
# (old row 152 -> new row 153, ... old row 263 -> new row 264).
#
# Inserting a full row at position 152 performs exactly that shift for
# every column, including carrying the date style (s="2") on column D
# down with it, so afterwards we only need to populate the brand-new
# row 152 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 152:263 down to 153:264 by inserting a new row.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A152").Value = 10
$ws.Range("B152").Value = "Vega Modelo de Temuco"
$ws.Range("C152").Value = "La Araucanía"
$ws.Range("D152").Value = 44762
$ws.Range("E152").Value = 9
$ws.Range("F152").Value = 100112039
$ws.Range("G152").Value = "Ciboulette"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 35
$ws.Range("K152").Value = 8000
$ws.Range("L152").Value = 8000
$ws.Range("M152").Value = 8000
$ws.Range("N152").Value = "$/docena de atados"
$ws.Range("O152").Value = "Provincia de Cautín"
$ws.Range("P152").Value = 2667
$ws.Range("Q152").Value = 3
$ws.Range("R152").Value = "Hortaliza"
